$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '29.772.32'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '1.621.57'
$ws.Range('E3').Value = '  +1.00%  '
Set-TextValue $ws.Range('D4') '0.991'
$ws.Range('E4').Value = '  -0.69%  '
Set-TextValue $ws.Range('D5') '212.72'
$ws.Range('E5').Value = '  +0.20%  '
Set-TextValue $ws.Range('D6') '0.521'
$ws.Range('E6').Value = '  +0.19%  '
Set-TextValue $ws.Range('D7') '0.989'
$ws.Range('E7').Value = '  -0.79%  '
Set-TextValue $ws.Range('D8') '29.37'
$ws.Range('E8').Value = '  +9.90%  '
$ws.Range('E9').Value = '  +3.61%  '
$ws.Range('E10').Value = '  +1.51%  '
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').Value = '1.848.29'
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('D13').Value = '1.638.29'
$ws.Range('E13').Value = '  +2.05%  '
Set-TextValue $ws.Range('D14') '0.567'
$ws.Range('E14').Value = '  +6.37%  '
$ws.Range('E15').Value = '  +5.70%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '29.770.66'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D17') '9.04'
$ws.Range('E17').Value = '  +18.13%  '
Set-TextValue $ws.Range('D18') '64.28'
$ws.Range('E18').Value = '  +1.91%  '
Set-TextValue $ws.Range('D19') '242.33'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('D20').Value = '0.0₃0710'
$ws.Range('E20').Value = '  +3.17%  '
Set-TextValue $ws.Range('D21') '0.992'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('E22').Value = '  +2.99%  '
Set-TextValue $ws.Range('D23') '9.68'
$ws.Range('E23').Value = '  +5.59%  '
$ws.Range('E24').Value = '  +1.19%  '
Set-TextValue $ws.Range('D25') '156.18'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('E26').Value = '  +2.66%  '
$ws.Range('E27').Value = '  +2.30%  '
$ws.Range('E28').Value = '  +3.53%  '
Set-TextValue $ws.Range('D29') '0.992'
$ws.Range('E29').Value = '  -0.56%  '
Set-TextValue $ws.Range('D30') '0.0487'
$ws.Range('E30').Value = '  +3.30%  '
$ws.Range('E31').Value = '  +2.89%  '
$ws.Range('E32').Value = '  +3.17%  '
$ws.Range('E33').Value = '  +3.44%  '
$ws.Range('D34').Value = '1.424.07'
$ws.Range('E34').Value = '  +0.80%  '
Set-TextValue $ws.Range('D35') '1.64'
$ws.Range('E35').Value = '  +7.34%  '
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('E37').Value = '  +2.06%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D38') '0.0171'
$ws.Range('E38').Value = '  +2.99%  '
$ws.Range('B39').Value = 'HuobiToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D39') '2.28'
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('E40').Value = '  +4.23%  '
$ws.Range('E41').Value = '  +3.48%  '
$ws.Range('E42').Value = '  +4.35%  '
$ws.Range('E43').Value = '  -0.02%  '
Set-TextValue $ws.Range('D44') '70.11'
$ws.Range('E44').Value = '  +6.84%  '
Set-TextValue $ws.Range('D45') '53.53'
$ws.Range('E45').Value = '  +1.24%  '
Set-TextValue $ws.Range('D46') '0.990'
$ws.Range('E46').Value = '  -0.76%  '
$ws.Range('E47').Value = '  +17.65%  '
Set-TextValue $ws.Range('D48') '5.46'
$ws.Range('E48').Value = '  +3.82%  '
$ws.Range('D49').Value = '1.760.34'
$ws.Range('E49').Value = '  +0.83%  '
Set-TextValue $ws.Range('D50') '88.00'
$ws.Range('E50').Value = '  +1.57%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0106'
$ws.Range('E51').Value = '  +2.11%  '
